$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "overall trial"
$ws.Range("D2").Value = "[1, 2]"
$ws.Range("E2").Value = "[(1, 2), (3, 4)]"
